# Insert a new data row at row 608 (pushing the existing rows 608.. down by one)
# and populate it with the new "Ajo" (Chino / Primera) price record described
# in the commit. All other rows keep their original values - they are simply
# shifted down by the insert, which Excel handles natively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(608).Insert()

$ws.Cells.Item(608, 1).Value  = 10
$ws.Cells.Item(608, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(608, 3).Value  = "La Araucanía"
$ws.Cells.Item(608, 4).Value  = 44714
$ws.Cells.Item(608, 5).Value  = 9
$ws.Cells.Item(608, 6).Value  = 100112003
$ws.Cells.Item(608, 7).Value  = "Ajo"
$ws.Cells.Item(608, 8).Value  = "Chino"
$ws.Cells.Item(608, 9).Value  = "Primera"
$ws.Cells.Item(608, 10).Value = 258
$ws.Cells.Item(608, 11).Value = 22000
$ws.Cells.Item(608, 12).Value = 22000
$ws.Cells.Item(608, 13).Value = 22000
$ws.Cells.Item(608, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(608, 15).Value = "China"
$ws.Cells.Item(608, 16).Value = 2200
$ws.Cells.Item(608, 17).Value = 10
$ws.Cells.Item(608, 18).Value = "Hortaliza"
